$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 131.57143
$ws.Range("I33").Value = 131.65
$ws.Range("K33").Value = 131.65
$ws.Range("M33").Value = 97.34999999999999

$ws.Range("H74").Value = 7357399
$ws.Range("J74").Value = 13893444
$ws.Range("L74").Value = 13893444
$ws.Range("N74").Value = -13895316

$ws.Range("H76").Value = 6947045.5
$ws.Range("I76").Value = 2979.5
$ws.Range("J76").Value = 9261734
$ws.Range("K76").Value = 2979.5
$ws.Range("L76").Value = 9261734
$ws.Range("M76").Value = -2664.5
$ws.Range("N76").Value = -9262364

$ws.Range("H77").Value = 7357399
$ws.Range("J77").Value = 13893444
$ws.Range("L77").Value = 69467220
$ws.Range("N77").Value = -69476580

$ws.Range("H79").Value = 6947045.5
$ws.Range("I79").Value = 2979.5
$ws.Range("J79").Value = 9261734
$ws.Range("K79").Value = 2979.5
$ws.Range("L79").Value = 9261734
$ws.Range("M79").Value = -1887.5
$ws.Range("N79").Value = -9263918

$ws.Range("H116").Value = 14709676
$ws.Range("I116").Value = 31251434
$ws.Range("K116").Value = 31251434
$ws.Range("M116").Value = -31247992

$ws.Range("H132").Value = 2199.6743
$ws.Range("I132").Value = 2199.6743
$ws.Range("K132").Value = 6599.0229
$ws.Range("M132").Value = -4069.0229

$ws.Range("H141").Value = 1685.6333
$ws.Range("I141").Value = 1394.16
$ws.Range("K141").Value = 4182.48
$ws.Range("M141").Value = 997.5199999999995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9117.769
$ws.Range("I32").Value = 6184.746
$ws.Range("J32").Value = 22218.6
$ws.Range("K32").Value = 6184.746
$ws.Range("L32").Value = 22218.6
$ws.Range("M32").Value = -5897.746
$ws.Range("N32").Value = -22792.6

$ws.Range("H41").Value = 2135.1
$ws.Range("I41").Value = 2356.375
$ws.Range("J41").Value = 1250
$ws.Range("K41").Value = 2356.375
$ws.Range("L41").Value = 1250
$ws.Range("M41").Value = -1942.375
$ws.Range("N41").Value = -2078

$ws.Range("H45").Value = 2373.9678
$ws.Range("I45").Value = 2245.0557
$ws.Range("J45").Value = 2552.4614
$ws.Range("K45").Value = 2245.0557
$ws.Range("L45").Value = 2552.4614
$ws.Range("M45").Value = -1868.0557
$ws.Range("N45").Value = -3306.4614

$ws.Range("H122").Value = 2329.2917
$ws.Range("I122").Value = 2267.762
$ws.Range("K122").Value = 6803.286
$ws.Range("M122").Value = -4353.286

$ws.Range("H132").Value = 9739.096
$ws.Range("I132").Value = 1617.0889
$ws.Range("J132").Value = 30044.111
$ws.Range("K132").Value = 4851.2667
$ws.Range("L132").Value = 90132.333
$ws.Range("M132").Value = -2321.2667
$ws.Range("N132").Value = -95192.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 59785
$ws.Range("J132").Value = 59785
$ws.Range("L132").Value = 59785
$ws.Range("N132").Value = -69905

$ws.Range("H134").Value = 3203.2341
$ws.Range("I134").Value = 3164.6086
$ws.Range("K134").Value = 9493.825800000001
$ws.Range("M134").Value = -6958.825800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4362.5
$ws.Range("I62").Value = 4380
$ws.Range("J62").Value = 4333.3335
$ws.Range("K62").Value = 4380
$ws.Range("L62").Value = 4333.3335
$ws.Range("M62").Value = -3756
$ws.Range("N62").Value = -5581.3335

$ws.Range("H65").Value = 4362.5
$ws.Range("I65").Value = 4380
$ws.Range("J65").Value = 4333.3335
$ws.Range("K65").Value = 21900
$ws.Range("L65").Value = 21666.6675
$ws.Range("M65").Value = -18780
$ws.Range("N65").Value = -27906.6675

$ws.Range("H99").Value = 20837066
$ws.Range("I99").Value = 3210
$ws.Range("K99").Value = 3210
$ws.Range("M99").Value = -1712

$ws.Range("H126").Value = 20837066
$ws.Range("I126").Value = 3210
$ws.Range("K126").Value = 9630
$ws.Range("M126").Value = -7160

$ws.Range("H132").Value = 2298.1428
$ws.Range("I132").Value = 1673.84
$ws.Range("J132").Value = 3858.9
$ws.Range("K132").Value = 5021.52
$ws.Range("L132").Value = 11576.7
$ws.Range("M132").Value = -2491.52
$ws.Range("N132").Value = -16636.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1240.425
$ws.Range("I5").Value = 982.0345
$ws.Range("J5").Value = 1921.6364
$ws.Range("K5").Value = 2946.1035
$ws.Range("L5").Value = 5764.9092
$ws.Range("M5").Value = -2834.1035
$ws.Range("N5").Value = -5988.9092

$ws.Range("H12").Value = 101.1
$ws.Range("J12").Value = 140.71428
$ws.Range("L12").Value = 422.14284
$ws.Range("N12").Value = -768.14284

$ws.Range("H92").Value = 41667300
$ws.Range("J92").Value = 1500
$ws.Range("L92").Value = 4500
$ws.Range("N92").Value = -6996

$ws.Range("H107").Value = 7116.7856
$ws.Range("I107").Value = 9313.429
$ws.Range("J107").Value = 526.8570999999999
$ws.Range("K107").Value = 27940.287
$ws.Range("L107").Value = 1580.5713
$ws.Range("M107").Value = -26020.287
$ws.Range("N107").Value = -5420.5713

$ws.Range("H113").Value = 924.53845
$ws.Range("J113").Value = 943.25
$ws.Range("L113").Value = 2829.75
$ws.Range("N113").Value = -7169.75

$ws.Range("H121").Value = 12821813
$ws.Range("J121").Value = 15152997
$ws.Range("L121").Value = 45458991
$ws.Range("N121").Value = -45461611

$ws.Range("H122").Value = 873.05884
$ws.Range("J122").Value = 1026.3077
$ws.Range("L122").Value = 9236.7693
$ws.Range("N122").Value = -14136.7693

$ws.Range("H131").Value = 643.7071
$ws.Range("J131").Value = 768.875
$ws.Range("L131").Value = 2306.625
$ws.Range("N131").Value = -12386.625

$ws.Range("H132").Value = 1213.5
$ws.Range("I132").Value = 1284.6666
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 11561.9994
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -9031.999400000001
$ws.Range("N132").Value = -14060

$ws.Range("H135").Value = 1240.425
$ws.Range("I135").Value = 982.0345
$ws.Range("J135").Value = 1921.6364
$ws.Range("K135").Value = 8838.3105
$ws.Range("L135").Value = 17294.7276
$ws.Range("M135").Value = -6303.3105
$ws.Range("N135").Value = -22364.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3686527
$ws.Range("I70").Value = 4866.6665
$ws.Range("J70").Value = 5694705.5
$ws.Range("K70").Value = 4866.6665
$ws.Range("L70").Value = 5694705.5
$ws.Range("M70").Value = -4596.6665
$ws.Range("N70").Value = -5695245.5

$ws.Range("H73").Value = 3686527
$ws.Range("I73").Value = 4866.6665
$ws.Range("J73").Value = 5694705.5
$ws.Range("K73").Value = 4866.6665
$ws.Range("L73").Value = 5694705.5
$ws.Range("M73").Value = -3930.6665
$ws.Range("N73").Value = -5696577.5

$ws.Range("H97").Value = 2749.6843
$ws.Range("I97").Value = 2596.5
$ws.Range("K97").Value = 2596.5
$ws.Range("M97").Value = -2100.5

$ws.Range("H122").Value = 6570.643
$ws.Range("I122").Value = 7061.25
$ws.Range("K122").Value = 21183.75
$ws.Range("M122").Value = -18733.75

$ws.Range("H132").Value = 12342.277
$ws.Range("I132").Value = 3316.8064
$ws.Range("K132").Value = 9950.4192
$ws.Range("M132").Value = -7420.4192

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 248236.8
$ws.Range("I132").Value = 346034.25
$ws.Range("J132").Value = 3743.1428
$ws.Range("K132").Value = 1038102.75
$ws.Range("L132").Value = 11229.4284
$ws.Range("M132").Value = -1035572.75
$ws.Range("N132").Value = -16289.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1352396.1
$ws.Range("J113").Value = 3378915
$ws.Range("L113").Value = 10136745
$ws.Range("N113").Value = -10141085
